$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: A16 gets the highlighted style (same as A2/A4/A12), B16 gets new note text
$ws.Range("A16").Interior.Color = 65535
$ws.Range("B16").Value = "My Experience - "

# Row 17: continuation note
$ws.Range("B17").Value = "Any type of method or property can be STUBBED / SHIMMED using MICROSOFT FAKES framework"

# Row 18: continuation note
$ws.Range("B18").Value = "You can STUB/SHIM that method (or) property in which it is defined.   No inheritance concepts observed here."

# Update the active selection to reflect the new last cell
$ws.Range("B18").Select()
